$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1935
$ws.Range("I2").Value = 4948
$ws.Range("J2").Value = 20682
$ws.Range("K2").Value = 102
$ws.Range("L2").Value = 5782
$ws.Range("M2").Value = 354
$ws.Range("N2").Value = 3595
$ws.Range("O2").Value = 13
$ws.Range("P2").Value = 71
$ws.Range("Q2").Value = 29
$ws.Range("R2").Value = 292
$ws.Range("S2").Value = 2264
$ws.Range("T2").Value = 3625
$ws.Range("U2").Value = 281
$ws.Range("V2").Value = 32004
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 32167
$ws.Range("Y2").Value = 62
$ws.Range("Z2").Value = 466
$ws.Range("AA2").Value = 196
